$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 370575.2
$ws.Range("J17").Value = 370575.2
$ws.Range("L17").Value = 1111725.6
$ws.Range("N17").Value = -1112061.6
$ws.Range("H112").Value = 24391906
$ws.Range("J112").Value = 1873.0857
$ws.Range("L112").Value = 5619.257100000001
$ws.Range("N112").Value = -7835.257100000001
$ws.Range("H129").Value = 876.89746
$ws.Range("I129").Value = 722.5294
$ws.Range("J129").Value = 996.1818
$ws.Range("K129").Value = 2167.5882
$ws.Range("L129").Value = 2988.5454
$ws.Range("M129").Value = 2832.4118
$ws.Range("N129").Value = -12988.5454
$ws.Range("H137").Value = 2223733.8
$ws.Range("I137").Value = 5264285.5
$ws.Range("J137").Value = 1792.1154
$ws.Range("K137").Value = 15792856.5
$ws.Range("L137").Value = 5376.3462
$ws.Range("M137").Value = -15790306.5
$ws.Range("N137").Value = -10476.3462
$ws.Range("H138").Value = 1145549.8
$ws.Range("I138").Value = 655.3182
$ws.Range("J138").Value = 2195036.2
$ws.Range("K138").Value = 1965.9546
$ws.Range("L138").Value = 6585108.600000001
$ws.Range("M138").Value = 3174.0454
$ws.Range("N138").Value = -6595388.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20998.564
$ws.Range("I32").Value = 23777.188
$ws.Range("J32").Value = 11471.857
$ws.Range("K32").Value = 23777.188
$ws.Range("L32").Value = 11471.857
$ws.Range("M32").Value = -23490.188
$ws.Range("N32").Value = -12045.857
$ws.Range("H61").Value = 43566330
$ws.Range("I61").Value = 52685264
$ws.Range("K61").Value = 52685264
$ws.Range("M61").Value = -52685052
$ws.Range("H74").Value = 14001207
$ws.Range("I74").Value = 16734662
$ws.Range("J74").Value = 333933
$ws.Range("K74").Value = 16734662
$ws.Range("L74").Value = 333933
$ws.Range("M74").Value = -16733788
$ws.Range("N74").Value = -335681
$ws.Range("H77").Value = 14001207
$ws.Range("I77").Value = 16734662
$ws.Range("J77").Value = 333933
$ws.Range("K77").Value = 83673310
$ws.Range("L77").Value = 1669665
$ws.Range("M77").Value = -83668942
$ws.Range("N77").Value = -1678401
$ws.Range("H132").Value = 45516.586
$ws.Range("I132").Value = 35172.367
$ws.Range("J132").Value = 64912
$ws.Range("K132").Value = 105517.101
$ws.Range("L132").Value = 194736
$ws.Range("M132").Value = -102987.101
$ws.Range("N132").Value = -199796
$ws.Range("H136").Value = 43566330
$ws.Range("I136").Value = 52685264
$ws.Range("K136").Value = 158055792
$ws.Range("M136").Value = -158053242

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1940.8572
$ws.Range("I134").Value = 2214.5715
$ws.Range("J134").Value = 1393.4286
$ws.Range("K134").Value = 6643.7145
$ws.Range("L134").Value = 4180.2858
$ws.Range("M134").Value = -4108.7145
$ws.Range("N134").Value = -9250.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3291.9697
$ws.Range("I31").Value = 1781.1052
$ws.Range("J31").Value = 5342.4287
$ws.Range("K31").Value = 1781.1052
$ws.Range("L31").Value = 5342.4287
$ws.Range("M31").Value = -1486.1052
$ws.Range("N31").Value = -5932.4287
$ws.Range("H34").Value = 3291.9697
$ws.Range("I34").Value = 1781.1052
$ws.Range("J34").Value = 5342.4287
$ws.Range("K34").Value = 1781.1052
$ws.Range("L34").Value = 5342.4287
$ws.Range("M34").Value = -1579.1052
$ws.Range("N34").Value = -5746.4287
$ws.Range("H58").Value = 30305090
$ws.Range("I58").Value = 52632944
$ws.Range("J58").Value = 2998.8572
$ws.Range("K58").Value = 52632944
$ws.Range("L58").Value = 2998.8572
$ws.Range("M58").Value = -52632741
$ws.Range("N58").Value = -3404.8572
$ws.Range("H132").Value = 36730.55
$ws.Range("I132").Value = 2300.8235
$ws.Range("J132").Value = 85506
$ws.Range("K132").Value = 6902.470499999999
$ws.Range("L132").Value = 256518
$ws.Range("M132").Value = -4372.470499999999
$ws.Range("N132").Value = -261578
$ws.Range("H134").Value = 32471.03
$ws.Range("I134").Value = 1907.3478
$ws.Range("J134").Value = 91051.414
$ws.Range("K134").Value = 5722.0434
$ws.Range("L134").Value = 273154.242
$ws.Range("M134").Value = -3187.0434
$ws.Range("N134").Value = -278224.242
$ws.Range("H136").Value = 30305090
$ws.Range("I136").Value = 52632944
$ws.Range("J136").Value = 2998.8572
$ws.Range("K136").Value = 157898832
$ws.Range("L136").Value = 8996.571599999999
$ws.Range("M136").Value = -157896282
$ws.Range("N136").Value = -14096.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 159.21739
$ws.Range("J38").Value = 208.8125
$ws.Range("L38").Value = 626.4375
$ws.Range("N38").Value = -1320.4375
$ws.Range("H113").Value = 539.30304
$ws.Range("I113").Value = 496.53333
$ws.Range("J113").Value = 574.94446
$ws.Range("K113").Value = 1489.59999
$ws.Range("L113").Value = 1724.83338
$ws.Range("M113").Value = 680.4000100000001
$ws.Range("N113").Value = -6064.83338
$ws.Range("H132").Value = 1286.3
$ws.Range("J132").Value = 1928.5714
$ws.Range("L132").Value = 17357.1426
$ws.Range("N132").Value = -22417.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1573.2646
$ws.Range("I102").Value = 1473.6786
$ws.Range("J102").Value = 2038
$ws.Range("K102").Value = 1473.6786
$ws.Range("L102").Value = 2038
$ws.Range("M102").Value = 148.3214
$ws.Range("N102").Value = -5282
$ws.Range("H132").Value = 127392.875
$ws.Range("I132").Value = 200900
$ws.Range("J132").Value = 93980.55
$ws.Range("K132").Value = 602700
$ws.Range("L132").Value = 281941.65
$ws.Range("M132").Value = -600170
$ws.Range("N132").Value = -287001.65

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3380.4443
$ws.Range("I40").Value = 2904
$ws.Range("J40").Value = 4333.3335
$ws.Range("K40").Value = 2904
$ws.Range("L40").Value = 4333.3335
$ws.Range("M40").Value = -2768
$ws.Range("N40").Value = -4605.3335
$ws.Range("H132").Value = 41327.81
$ws.Range("I132").Value = 1473
$ws.Range("J132").Value = 81182.62
$ws.Range("K132").Value = 4419
$ws.Range("L132").Value = 243547.86
$ws.Range("M132").Value = -1889
$ws.Range("N132").Value = -248607.86
$ws.Range("H136").Value = 201350.4
$ws.Range("I136").Value = 334368
$ws.Range("J136").Value = 144342.86
$ws.Range("K136").Value = 1003104
$ws.Range("L136").Value = 433028.58
$ws.Range("M136").Value = -1000554
$ws.Range("N136").Value = -438128.58

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 65756.71000000001
$ws.Range("I132").Value = 39403.848
$ws.Range("K132").Value = 118211.544
$ws.Range("M132").Value = -115681.544
$ws.Range("H136").Value = 65731.836
$ws.Range("I136").Value = 39738
$ws.Range("J136").Value = 200899.8
$ws.Range("K136").Value = 119214
$ws.Range("L136").Value = 602699.3999999999
$ws.Range("M136").Value = -116664
$ws.Range("N136").Value = -607799.3999999999

Write-Host "Applied all Hades_Profits updates"